# Applies the "Log and Excel Updated" edit:
#  - Job sheet: add column C ("Job Title Already exit") for rows 1-2
#  - SalComp sheet: change A1 "Food All2" -> "Food All3"; add B1 "Already exit"
#  - Selections / active sheet updated to match the authored state

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsJob = $wb.Worksheets.Item("Job")
$wsSalComp = $wb.Worksheets.Item("SalComp")

# --- Job sheet: new column C with "Job Title Already exit" ---
$wsJob.Range("C1").Value = "Job Title Already exit"
$wsJob.Range("C2").Value = "Job Title Already exit"

# --- SalComp sheet: update A1 text and add B1 ---
$wsSalComp.Range("A1").Value = "Food All3"
$wsSalComp.Range("B1").Value = "Already exit"

# --- selections to match the target state ---
$wsJob.Range("C17").Select() | Out-Null
$wsSalComp.Range("F7").Select() | Out-Null
$wsLogin.Range("C16").Select() | Out-Null

# --- make Login the active sheet/tab ---
$wsLogin.Activate() | Out-Null
